$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "HK_G_acc_SD"
$ws.Cells.Item(2, 1).Value = 64.769647696476966
$ws.Cells.Item(3, 1).Value = 64.769647696476966
$ws.Cells.Item(4, 1).Value = 65.040650406504056
$ws.Cells.Item(5, 1).Value = 64.498644986449861
$ws.Cells.Item(6, 1).Value = 63.956639566395665
$ws.Cells.Item(7, 1).Value = 64.498644986449861
$ws.Cells.Item(8, 1).Value = 65.853658536585371
$ws.Cells.Item(9, 1).Value = 68.834688346883468
$ws.Cells.Item(10, 1).Value = 69.647696476964768
$ws.Cells.Item(11, 1).Value = 67.208672086720867
$ws.Cells.Item(12, 1).Value = 65.311653116531161
$ws.Cells.Item(13, 1).Value = 68.563685636856363
$ws.Cells.Item(14, 1).Value = 69.105691056910572
$ws.Cells.Item(15, 1).Value = 68.834688346883468
$ws.Cells.Item(16, 1).Value = 68.834688346883468
$ws.Cells.Item(17, 1).Value = 67.208672086720867
$ws.Cells.Item(18, 1).Value = 67.750677506775077
$ws.Cells.Item(19, 1).Value = 71.002710027100264
$ws.Cells.Item(20, 1).Value = 66.937669376693762
$ws.Cells.Item(21, 1).Value = 66.124661246612476
$ws.Cells.Item(22, 1).Value = 65.853658536585371
$ws.Cells.Item(23, 1).Value = 59.078590785907856
$ws.Cells.Item(24, 1).Value = 59.078590785907856
$ws.Cells.Item(25, 1).Value = 57.72357723577236
$ws.Cells.Item(26, 1).Value = 67.750677506775077
$ws.Cells.Item(27, 1).Value = 67.208672086720867
$ws.Cells.Item(28, 1).Value = 69.918699186991873
$ws.Cells.Item(29, 1).Value = 68.563685636856363
$ws.Cells.Item(30, 1).Value = 68.021680216802167
$ws.Cells.Item(31, 1).Value = 69.105691056910572
$ws.Cells.Item(32, 1).Value = 59.078590785907856
$ws.Cells.Item(33, 1).Value = 59.891598915989164
$ws.Cells.Item(34, 1).Value = 60.433604336043359
$ws.Cells.Item(35, 1).Value = 57.452574525745263
$ws.Cells.Item(36, 1).Value = 55.826558265582662
$ws.Cells.Item(37, 1).Value = 60.975609756097562
$ws.Cells.Item(38, 1).Value = 56.36856368563685
$ws.Cells.Item(39, 1).Value = 56.36856368563685
$ws.Cells.Item(40, 1).Value = 56.36856368563685
$ws.Cells.Item(41, 1).Value = 67.750677506775077
$ws.Cells.Item(42, 1).Value = 68.834688346883468
$ws.Cells.Item(43, 1).Value = 69.376693766937663
$ws.Cells.Item(44, 1).Value = 66.937669376693762
$ws.Cells.Item(45, 1).Value = 65.311653116531161
$ws.Cells.Item(46, 1).Value = 67.208672086720867
$ws.Cells.Item(47, 1).Value = 65.040650406504056
$ws.Cells.Item(48, 1).Value = 66.395663956639567
$ws.Cells.Item(49, 1).Value = 67.750677506775077
